# Model_Parameters.xlsx edit: change bound formulas from /10,*10 to /5,*5
# for the "fitted" tunable-parameter rows, and preserve the original bound
# values in two new columns ("lower bound orig" / "upper bound orig").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Add the two new header cells O1 / P1 (copy format from N1, the
#    existing "logscale fitting" header) and set their text.
# ---------------------------------------------------------------------
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("O1").Value = "lower bound orig"
$ws.Range("P1").Value = "upper bound orig"

# ---------------------------------------------------------------------
# 2) Copy the CURRENT (pre-edit) lower/upper bound values into the new
#    O / P columns as static values, for every row that currently has a
#    lower/upper bound in columns L / M.
# ---------------------------------------------------------------------
$boundRows = @(2,3,4,5,8,10,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44)
foreach ($r in $boundRows) {
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 12).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($r, 13).Value2
}

# ---------------------------------------------------------------------
# 3) Change the lower/upper bound formulas for the rows whose bound was
#    a simple /10 , *10 of the full-model value to instead use /5 , *5.
#    (Row 2's upper bound formula, =E2*100, is left exactly as-is - that
#    inconsistency is present in the target workbook too.)
# ---------------------------------------------------------------------
$ws.Range("L2").Formula = "=E2/5"

$ws.Range("L3").Formula = "=E3/5"
$ws.Range("M3").Formula = "=E3*5"

$ws.Range("L8").Formula = "=E8/5"
$ws.Range("M8").Formula = "=E8*5"

$ws.Range("L10").Formula = "=E10/5"
$ws.Range("M10").Formula = "=E10*5"

$ws.Range("L12").Formula = "=E12/5"
$ws.Range("M12").Formula = "=E12*5"

$ws.Range("L13:L26").Formula = "=E13/5"
$ws.Range("M13:M26").Formula = "=E13*5"

# Writing a fresh .Formula to a previously-unstyled cell can pick up a
# stray number-format style; put those cells back to the default
# "Normal" style so they stay unstyled, matching the original workbook.
$ws.Range("L12:M26").Style = "Normal"

# ---------------------------------------------------------------------
# 4) Misc cosmetic bits captured by the diff: active cell selection.
# ---------------------------------------------------------------------
$null = $ws.Range("J34").Select()
